$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "Upcoming" announcements text box on slide 1.
$sh = $s.Shapes.Item("TextBox 3")
$tr = $sh.TextFrame.TextRange

# Paragraph 3 is currently "Lab 4 released 6:00 pm Thursday." — insert the new
# "Online Term Test Review" bullet immediately before it (it inherits the
# bullet / paragraph formatting of paragraph 3 automatically).
$para3 = $tr.Paragraphs(3)
$newPara = $para3.InsertBefore("Online Term Test Review, Tuesday at 7:30 pm.`r")

# Colour the whole new line green, matching the sibling bullets.
$newPara.Font.Color.RGB = 10092390   # 66FF99 (R=0x66,G=0xFF,B=0x99 -> R+G*256+B*65536)

# The trailing period should be purple like the other bullets' periods.
# newPara.Length includes the paragraph-mark, so the period sits two
# characters before the end of the range.
$periodPos = $newPara.Start + $newPara.Length - 2
$period = $tr.Characters($periodPos, 1)
$period.Font.Color.RGB = 16751052    # CC99FF (R=0xCC,G=0x99,B=0xFF -> R+G*256+B*65536)
